$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "plate" column (B) data values: numeric plate numbers -> text plate labels ---
$ws.Range("B2:B25").Value  = "Plate1"
$ws.Range("B26:B49").Value = "Plate2"
$ws.Range("B50:B73").Value = "Plate3"

# --- Column A no longer carries the bespoke "date" number-format style; drop it back
#     to the workbook default (this also clears the per-cell s="3" on A2:A73, and the
#     col-level default style, and slightly narrows the column). ---
$ws.Columns(1).ClearFormats()
$ws.Columns(1).ColumnWidth = 9.666666666666666
$ws.Columns(4).ColumnWidth = 11.498697916666666
$ws.Columns(6).ColumnWidth = 11.666666666666666

# --- Re-apply the shared header look (9pt Helvetica, thin grey border) to A1, matching
#     the rest of row 1 (B1:F1) now that ClearFormats() wiped it. ---
$ws.Range("A1").Font.Name = "Helvetica"
$ws.Range("A1").Font.Size = 9
$ws.Range("A1").Borders.LineStyle = 1
$ws.Range("A1").Borders.Color = 12566463

# --- Move the active selection/view as it was left after the edit ---
$ws.Range("D79").Select()
